# v0.8: improved excel loading. fixed missing pi and other fields
#
# The original workbook has a single sheet "data" holding a flat
# report/field metadata table. This edit:
#   1. Inserts a new blank "Sheet1" in front of the existing sheet.
#   2. Renames the original "data" sheet to "ADIDO Metadata" and makes
#      it the active sheet.
#   3. Re-shapes its table: new headers (File Name / Field Name /
#      Business Description / Classification / PCI / PI / Data
#      Treatment), a new constant "Business Description" column, the
#      PI/PCI columns swapped, and the previously-skipped report_2/name
#      row filled back in (closing the old row-5 gap).

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new blank worksheet before the existing one --------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet1"

# --- 2. Rename the original data sheet -------------------------------------
$meta = $wb.Worksheets.Item("data")
$meta.Name = "ADIDO Metadata"

# --- 3. Rebuild the metadata table -----------------------------------------
$meta.UsedRange.ClearContents()

# Main block (File Name .. PI) pasted in first ...
$headersAF = @("File Name", "Field Name", "Business Description", "Classification", "PCI", "PI")
for ($col = 1; $col -le $headersAF.Length; $col++) {
    $meta.Cells.Item(1, $col).Value = $headersAF[$col - 1]
}

$rowsAF = @(
    @("report_1", "id", "not available", "internal", $false, $false),
    @("report_1", "number", "not available", "internal", $false, $false),
    @("report_1", "credit_card_number", "not available", "internal", $false, $false),
    @("report_2", "id", "not available", "internal", $false, $false),
    @("report_2", "name", "not available", "internal", $false, $false),
    @("report_2", "full_name", "not available", "internal", $false, $false),
    @("report_2", "date_of_birth", "not available", "internal", $false, $false),
    @("report_2", "credit_score", "not available", "internal", $false, $false)
)

for ($i = 0; $i -lt $rowsAF.Length; $i++) {
    $r = $i + 2
    $row = $rowsAF[$i]
    for ($col = 1; $col -le $row.Length; $col++) {
        $meta.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

# ... the "Data Treatment" column is added afterwards.
$meta.Cells.Item(1, 7).Value = "Data Treatment"

$colG = @("default", "default", "masked", "default", "default", "default", "masked", "rounded")
for ($i = 0; $i -lt $colG.Length; $i++) {
    $meta.Cells.Item($i + 2, 7).Value = $colG[$i]
}

# --- 4. Page setup / selection / active sheet -------------------------------
$meta.PageSetup.Orientation = 1

$meta.Range("G2").Select()
$meta.Activate()
